# Updated DB connections to use DATABASE_URL for Render
#
# Adds the new Credits/Cost comparison rows (8 & 9) on the JobSeekers sheet
# and the two percentage "growth" formulas in row 10 (K10, L10), then
# updates the sheet's on-screen selection to match the author's final
# cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JobSeekers")
$ws.Activate()

# --- New data: rows 8 and 9 (columns K and L) ---
$ws.Cells.Item(8, 11).Value = 1400   # K8
$ws.Cells.Item(8, 12).Value = 114    # L8

$ws.Cells.Item(9, 11).Value = 247    # K9
$ws.Cells.Item(9, 12).Value = 70     # L9

# --- Row 10: new percentage formulas next to the existing B10 rate ---
$ws.Cells.Item(10, 11).Formula = "=K9/K8"          # K10
$ws.Cells.Item(10, 11).Style = "Percent"

$ws.Cells.Item(10, 12).Formula = "=(L9-L8)/L8"     # L10
$ws.Cells.Item(10, 12).Style = "Percent"

# --- View state: selection moves to M10 (scrolled so column F is leftmost) ---
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("M10").Select() | Out-Null

$wb.Saved = $false
